$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "University Of Washington FullStack Developer BootCamp" paragraph:
#    split " Washington " + "FullStack" + " Developer " + "BootCamp" (which
#    carried spellStart/spellEnd proofErr flags around the two run-together
#    words) into properly spaced runs: "Full Stack" and "Boot Camp", and
#    drop the now-stale spell-check markers.
# ---------------------------------------------------------------------------
$rFull = $d.Content
$null = $rFull.Find.Execute("FullStack")
$eduPara = $rFull.Paragraphs(1)

$eduXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
  '<w:body>' +
  '<w:p w14:paraId="3A1909AB" w14:textId="76719B01" w:rsidR="001C45F8" w:rsidRDefault="001C45F8">' +
    '<w:pPr>' +
      '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
      '<w:ind w:left="720"/>' +
      '<w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve">University </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>Of</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> Washington Full</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>Stack Developer Boot</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>Camp</w:t></w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$null = $eduPara.Range.InsertXML($eduXml)

# ---------------------------------------------------------------------------
# 2) "September, 2020 to Present" paragraph: collapse the " 20" / "20" /
#    " to " / "Present" runs into a single " 2020 to Present" run (keeping
#    the "September," gramStart/gramEnd markers untouched).
# ---------------------------------------------------------------------------
$rWork = $d.Content
$null = $rWork.Find.Execute("2020 to Present")
$workPara = $rWork.Paragraphs(1)

$workXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
  '<w:body>' +
  '<w:p w14:paraId="679C1D12" w14:textId="5E60FF38" w:rsidR="001C45F8" w:rsidRDefault="001C45F8" w:rsidP="001C45F8">' +
    '<w:pPr>' +
      '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
      '<w:ind w:left="720"/>' +
      '<w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>' +
    '</w:pPr>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>September,</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> 2020 to Present</w:t></w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$null = $workPara.Range.InsertXML($workXml)
